$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 104.54
$ws.Range("E2").Value = 54.9
$ws.Range("F2").Value = 14.4
$ws.Range("N2").Value = 50.68470204858703

$ws.Range("D3").Value = 22.04
$ws.Range("E3").Value = 48.7
$ws.Range("F3").Value = 10.2
$ws.Range("N3").Value = 50.68470204858703
